$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("B2:B61").NumberFormat = "@"

$ws.Cells.Item(2,1).Value = 71
$ws.Cells.Item(2,2).Value = "81"
$ws.Cells.Item(3,1).Value = 72
$ws.Cells.Item(3,2).Value = "64"
$ws.Cells.Item(4,1).Value = 74
$ws.Cells.Item(4,2).Value = "66"
$ws.Cells.Item(5,1).Value = 75
$ws.Cells.Item(5,2).Value = "26"
$ws.Cells.Item(6,1).Value = 76
$ws.Cells.Item(6,2).Value = "36"
$ws.Cells.Item(7,1).Value = 77
$ws.Cells.Item(7,2).Value = "16"
$ws.Cells.Item(8,1).Value = 78
$ws.Cells.Item(8,2).Value = "76"
$ws.Cells.Item(9,1).Value = 79
$ws.Cells.Item(9,2).Value = "24"
$ws.Cells.Item(10,1).Value = 80
$ws.Cells.Item(10,2).Value = "80"
$ws.Cells.Item(11,1).Value = 81
$ws.Cells.Item(11,2).Value = "71"
$ws.Cells.Item(12,1).Value = 83
$ws.Cells.Item(12,2).Value = "84"
$ws.Cells.Item(13,1).Value = 84
$ws.Cells.Item(13,2).Value = "83"
$ws.Cells.Item(14,1).Value = 71
$ws.Cells.Item(14,2).Value = "56"
$ws.Cells.Item(15,1).Value = 72
$ws.Cells.Item(15,2).Value = "68"
$ws.Cells.Item(16,1).Value = 74
$ws.Cells.Item(16,2).Value = "54"
$ws.Cells.Item(17,1).Value = 75
$ws.Cells.Item(17,2).Value = "33"
$ws.Cells.Item(18,1).Value = 76
$ws.Cells.Item(18,2).Value = "78"
$ws.Cells.Item(19,1).Value = 77
$ws.Cells.Item(19,2).Value = "22"
$ws.Cells.Item(20,1).Value = 78
$ws.Cells.Item(20,2).Value = "50"
$ws.Cells.Item(21,1).Value = 79
$ws.Cells.Item(21,2).Value = "29"
$ws.Cells.Item(22,1).Value = 80
$ws.Cells.Item(22,2).Value = "67"
$ws.Cells.Item(23,1).Value = 81
$ws.Cells.Item(23,2).Value = "55"
$ws.Cells.Item(24,1).Value = 83
$ws.Cells.Item(24,2).Value = "6"
$ws.Cells.Item(25,1).Value = 84
$ws.Cells.Item(25,2).Value = "6"
$ws.Cells.Item(26,1).Value = 71
$ws.Cells.Item(26,2).Value = "55"
$ws.Cells.Item(27,1).Value = 72
$ws.Cells.Item(27,2).Value = "74"
$ws.Cells.Item(28,1).Value = 74
$ws.Cells.Item(28,2).Value = "68"
$ws.Cells.Item(29,1).Value = 75
$ws.Cells.Item(29,2).Value = "42"
$ws.Cells.Item(30,1).Value = 76
$ws.Cells.Item(30,2).Value = "79"
$ws.Cells.Item(31,1).Value = 77
$ws.Cells.Item(31,2).Value = "14"
$ws.Cells.Item(32,1).Value = 78
$ws.Cells.Item(32,2).Value = "49"
$ws.Cells.Item(33,1).Value = 79
$ws.Cells.Item(33,2).Value = "27"
$ws.Cells.Item(34,1).Value = 80
$ws.Cells.Item(34,2).Value = "50"
$ws.Cells.Item(35,1).Value = 81
$ws.Cells.Item(35,2).Value = "56"
$ws.Cells.Item(36,1).Value = 83
$ws.Cells.Item(36,2).Value = "8"
$ws.Cells.Item(37,1).Value = 84
$ws.Cells.Item(37,2).Value = "3"
$ws.Cells.Item(38,1).Value = 71
$ws.Cells.Item(38,2).Value = "63"
$ws.Cells.Item(39,1).Value = 72
$ws.Cells.Item(39,2).Value = "60"
$ws.Cells.Item(40,1).Value = 74
$ws.Cells.Item(40,2).Value = "72"
$ws.Cells.Item(41,1).Value = 75
$ws.Cells.Item(41,2).Value = "83"
$ws.Cells.Item(42,1).Value = 76
$ws.Cells.Item(42,2).Value = "24"
$ws.Cells.Item(43,1).Value = 77
$ws.Cells.Item(43,2).Value = "13"
$ws.Cells.Item(44,1).Value = 78
$ws.Cells.Item(44,2).Value = "44"
$ws.Cells.Item(45,1).Value = 79
$ws.Cells.Item(45,2).Value = "30"
$ws.Cells.Item(46,1).Value = 80
$ws.Cells.Item(46,2).Value = "51"
$ws.Cells.Item(47,1).Value = 81
$ws.Cells.Item(47,2).Value = "63"
$ws.Cells.Item(48,1).Value = 83
$ws.Cells.Item(48,2).Value = "3"
$ws.Cells.Item(49,1).Value = 84
$ws.Cells.Item(49,2).Value = "12"
$ws.Cells.Item(50,1).Value = 71
$ws.Cells.Item(50,2).Value = "41"
$ws.Cells.Item(51,1).Value = 72
$ws.Cells.Item(51,2).Value = "67"
$ws.Cells.Item(52,1).Value = 74
$ws.Cells.Item(52,2).Value = "60"
$ws.Cells.Item(53,1).Value = 75
$ws.Cells.Item(53,2).Value = "30"
$ws.Cells.Item(54,1).Value = 76
$ws.Cells.Item(54,2).Value = "46"
$ws.Cells.Item(55,1).Value = 77
$ws.Cells.Item(55,2).Value = "10"
$ws.Cells.Item(56,1).Value = 78
$ws.Cells.Item(56,2).Value = "36"
$ws.Cells.Item(57,1).Value = 79
$ws.Cells.Item(57,2).Value = "28"
$ws.Cells.Item(58,1).Value = 80
$ws.Cells.Item(58,2).Value = "63"
$ws.Cells.Item(59,1).Value = 81
$ws.Cells.Item(59,2).Value = "47"
$ws.Cells.Item(60,1).Value = 83
$ws.Cells.Item(60,2).Value = "12"
$ws.Cells.Item(61,1).Value = 84
$ws.Cells.Item(61,2).Value = "2"

$ws.Range("B2:B61").Style = "Normal"
